$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at G (shifts old G->H, H->I, I->J, J->K, K->L, N->O)
$ws.Columns("G:G").Insert()

# 2. Point the status-summary defined names at their new column (K -> L)
$wb.Names.Item("ln_completed").RefersTo = "=Hoja1!`$L`$10"
$wb.Names.Item("ln_in_progress").RefersTo = "=Hoja1!`$L`$9"
$wb.Names.Item("ln_pending").RefersTo = "=Hoja1!`$L`$8"

# 3. New "Github Actions" header + SI/NO values for every service row
$ws.Range("G3").Value = "Github Actions"
$ws.Range("G4").Value = "NO"
$ws.Range("G5").Value = "NO"
$ws.Range("G6").Value = "NO"
$ws.Range("G8").Value = "NO"
$ws.Range("G9").Value = "NO"
$ws.Range("G10").Value = "NO"
$ws.Range("G11").Value = "NO"
$ws.Range("G13").Value = "NO"
$ws.Range("G14").Value = "NO"
$ws.Range("G15").Value = "NO"
$ws.Range("G16").Value = "NO"
$ws.Range("G17").Value = "NO"
$ws.Range("G18").Value = "NO"

# 4. gateway-server (row 6) is now deployed to Azure with Github Actions still pending,
#    plus its container URL/port
$ws.Range("F6").Value = "SI"
$ws.Range("H6").Value = "http://gateway-service.eastus.azurecontainer.io/"
$ws.Range("I6").Value = 80

# 5. Conditional formatting: stretch the existing ranges to include the new column
$ws.Range("K7:K10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("L7:L10"))
$r1 = $ws.Range("E25:H1048576,E7:H20")
$cf1 = $ws.Range("K7:K10 E25:H1048576 E7:H20")

$ws.Range("N1:N5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("O1:O5,O8:O9"))

$ws.Range("E3:H3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E3:I3"))
$ws.Range("E5:H5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E5:I5"))
$ws.Range("F6:H6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("F6:I6"))
$ws.Range("F1:H1048576").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("F1:I1048576"))

# New conditional formatting for the two freshly-populated cells in row 6
$cfFG6 = $ws.Range("F6:G6").FormatConditions
$rFG1 = $cfFG6.Add(1, 3, '"COMPLETADO"')
$rFG1.Interior.Color = 13561798
$rFG2 = $cfFG6.Add(1, 3, '"EN PROGRESO"')
$rFG3 = $cfFG6.Add(1, 3, '"PENDIENTE"')

$cfH6 = $ws.Range("H6").FormatConditions
$rH1 = $cfH6.Add(1, 3, '"COMPLETADO"')
$rH2 = $cfH6.Add(1, 3, '"EN PROGRESO"')
$rH3 = $cfH6.Add(1, 3, '"PENDIENTE"')

# 6. Selection moves one column left (onto the Estado column) after the insert
$ws.Range("E5").Select()
